$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(87.15747049631544, 85.56960883729815, 86.74909249621672, 86.43756662660813)
    3 = @(96.89236841920233, 93.97978612548049, 97.74026796983379, 96.11139236723272)
    4 = @(98.99688545465644, 99.01798399007733, 98.62829296414554, 99.10502191216564)
    5 = @(98.69190515134034, 98.53923888870015, 98.57874915973579, 98.64413817508448)
    6 = @(98.42840706434657, 98.4513514054043,  98.51121505988981, 98.34774029915036)
    7 = @(97.84222378130558, 97.78371362265267, 98.09169748338373, 97.63473444963604)
    8 = @(97.90018380601698, 96.99225428792676, 97.21957300024354, 97.6788984342578)
    9 = @(96.04604249573356, 95.96238544459246, 96.5745243230481,  95.7771745445792)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}
